# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" values in E16/E17 are swapped (2107 <-> 2104) and the
# "Salario Basico" values in G16/G17 are updated from 908526 to 877803.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Swap the Periodo Mora values for the two detail rows.
$ws.Range("E16").Value = "2104"
$ws.Range("E17").Value = "2107"

# Update Salario Basico for both detail rows.
$ws.Range("G16").Value = 877803
$ws.Range("G17").Value = 877803
